# Linking between CRS and SIQ
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRS")

# --- Unmerge D2:D6 (keep E2:E6 merged) so D6 can hold its own value/style ---
$ws.Range("D2:D6").UnMerge()

# --- Row 6: add the "SIQ ID" header in column D, matching the other header cells ---
$ws.Range("A6").Copy()
$ws.Range("D6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D6").Value = "SIQ ID"

# --- Rows 2:5 column D now stand alone (no more merge) - centered vertically, wrap, no border, no horizontal centering ---
$dRange = $ws.Range("D2:D5")
$dRange.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none
$dRange.HorizontalAlignment = 1             # xlGeneral
$dRange.VerticalAlignment = -4108           # xlCenter
$dRange.WrapText = $true

# --- Row 8: update the requirement description text and link it to its SIQ ---
$ws.Range("C8").Value = "The system should be Initialized with the last temperature and the fan speed reading before the system powered off"
$ws.Range("D8").Value = "SIQ_5"

# --- Rows 10-14: link each CRS row to its corresponding SIQ id ---
$ws.Range("C14").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("D10").Value = "SIQ_1"
$ws.Range("D11").Value = "SIQ_4"
$ws.Range("D12").Value = "SIQ_3"
$ws.Range("D13").Value = "SIQ_6"
$ws.Range("D14").Value = "SIQ_2"

# --- restore the active view/selection ---
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D7").Select()
